$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.478.79"
$ws.Range("E2").Value2 = "  +3.76%  "

$ws.Range("D3").Value2 = "1.911.13"
$ws.Range("E3").Value2 = "  +2.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "333.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.79%  "

$ws.Range("E6").Value2 = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4674"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  +1.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.4101"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +2.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "47.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.08039"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +2.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "1.013"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +2.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "22.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +4.84%  "

$ws.Range("D13").Value2 = "1.921.30"
$ws.Range("E13").Value2 = "  +2.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.974"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +2.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "7.178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +2.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "89.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +1.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +0.11%  "

$ws.Range("E18").Value2 = "  +1.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.06586"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +0.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "17.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +3.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +0.30%  "

$ws.Range("D22").Value2 = "29.425.07"
$ws.Range("E22").Value2 = "  +3.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +4.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "11.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +5.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.215"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -1.50%  "

$ws.Range("D26").Value2 = "2.113.08"
$ws.Range("E26").Value2 = "  +1.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "154.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -1.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "19.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +2.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "5.768"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +8.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "117.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.066"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +11.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.09465"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +1.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.427"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +2.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "3.565"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.406"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +3.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.06116"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +1.41%  "

$ws.Range("E38").Value2 = "  +2.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "8.417"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +1.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.177"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +1.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.5887"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +2.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.1842"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +1.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "10.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +1.55%  "

$ws.Range("E44").Value2 = "  +0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "2.351"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +2.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.07507"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +5.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.5569"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +2.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "12.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +2.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.928"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "113.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +1.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.2960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +9.87%  "

